$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Status text: "In Translation" -> "Ready for handoff"
$wsOverview.Range("E2").Value() = "Ready for handoff"
$wsOverview.Range("F2").Value() = "Ready for handoff"
$wsZhCn.Range("C2").Value()     = "Ready for handoff"
$wsDeDe.Range("C2").Value()     = "Ready for handoff"

# Latest HO Xliff Generate Date / Latest Handoff Datetime
# Overview!G2 and de-de!H2 shared the same timestamp string before the edit
$wsOverview.Range("G2").Value() = "2016-08-29 18:42:40"
$wsDeDe.Range("H2").Value()     = "2016-08-29 18:42:40"

# zh-cn!H2 had its own distinct timestamp string
$wsZhCn.Range("H2").Value() = "2016-08-29 18:42:36"

# Column width adjustments (report regenerated -> wider Status/date columns)
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
